$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the value of C4 to "FEITO" (matches existing shared string used in C2, C3, C6, C8, C9)
$ws.Range("C4").Value = "FEITO"

# Update the active selection to C4, matching the saved selection state in the diff
$ws.Range("C4").Select()
